$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "68.921.92"; E = "  -0.12%  " }
    @{ Row = 3; D = "3.761.74"; E = "  -1.78%  " }
    @{ Row = 4; D = $null; E = "  +0.30%  " }
    @{ Row = 5; D = "'629.44"; E = "  +0.62%  " }
    @{ Row = 6; D = "'165.34"; E = "  +0.04%  " }
    @{ Row = 7; D = "3.759.40"; E = "  -1.70%  " }
    @{ Row = 8; D = $null; E = "  -0.06%  " }
    @{ Row = 9; D = "'0.520"; E = "  +0.24%  " }
    @{ Row = 10; D = $null; E = "  -2.08%  " }
    @{ Row = 11; D = $null; E = "  +0.65%  " }
    @{ Row = 12; D = "'6.77"; E = "  +1.46%  " }
    @{ Row = 13; D = $null; E = "  -4.24%  " }
    @{ Row = 14; D = "'34.80"; E = "  -2.64%  " }
    @{ Row = 15; D = "4.391.67"; E = "  -1.56%  " }
    @{ Row = 16; D = "3.760.51"; E = "  -1.25%  " }
    @{ Row = 17; D = "68.907.35"; E = "  -0.05%  " }
    @{ Row = 18; D = "'17.63"; E = "  -3.00%  " }
    @{ Row = 19; D = $null; E = "  +0.08%  " }
    @{ Row = 20; D = $null; E = "  -1.86%  " }
    @{ Row = 21; D = "'461.17"; E = "  -1.53%  " }
    @{ Row = 22; D = "'9.46"; E = "  -2.85%  " }
    @{ Row = 23; D = $null; E = "  -0.38%  " }
    @{ Row = 24; D = "'81.96"; E = "  -2.38%  " }
    @{ Row = 25; D = $null; E = "  -4.70%  " }
    @{ Row = 26; D = $null; E = "  +0.44%  " }
    @{ Row = 27; D = $null; E = "  -1.61%  " }
    @{ Row = 28; D = "'10.09"; E = "  +0.29%  " }
    @{ Row = 29; D = $null; E = "  -0.10%  " }
    @{ Row = 30; D = "3.909.50"; E = "  -1.50%  " }
    @{ Row = 31; D = "'2.27"; E = "  +0.95%  " }
    @{ Row = 32; D = "'2.67"; E = "  +0.12%  " }
    @{ Row = 33; D = "'7.04"; E = "  -3.30%  " }
    @{ Row = 34; D = $null; E = "  -3.25%  " }
    @{ Row = 35; D = "'0.175"; E = "  +18.12%  " }
    @{ Row = 36; D = "'0.999"; E = "  -0.04%  " }
    @{ Row = 37; D = "3.711.85"; E = "  -1.49%  " }
    @{ Row = 38; D = "'8.86"; E = "  -2.44%  " }
    @{ Row = 39; D = $null; E = "  -1.51%  " }
    @{ Row = 40; D = "'3.27"; E = "  +1.16%  " }
    @{ Row = 41; D = $null; E = "  -2.39%  " }
    @{ Row = 42; D = "'0.999"; E = "  +0.16%  " }
    @{ Row = 43; D = "'0.960"; E = "  -2.06%  " }
    @{ Row = 44; D = $null; E = "  +0.00%  " }
    @{ Row = 45; D = "'156.91"; E = "  -0.58%  " }
    @{ Row = 46; D = $null; E = "  +3.50%  " }
    @{ Row = 47; D = $null; E = "  +0.80%  " }
    @{ Row = 48; D = "'47.06"; E = "  +0.73%  " }
    @{ Row = 49; D = "'42.62"; E = "  -0.89%  " }
    @{ Row = 50; D = $null; E = "  -2.47%  " }
    @{ Row = 51; D = $null; E = "  -0.93%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
